# [fix][asset] Fix column bug to download asset
#
# - "part" sheet: rename the "硬盘"/"RAID" headers to "系统盘"/"数据盘",
#   and add two new trailing columns "模块" / "配件变更".
# - Widen the new/renamed data columns (D:K) on the "part" sheet to match
#   the existing wide columns.
# - Move the saved cursor/selection on "asset" to H20 and on "part" to A1.

$wb = $excel.ActiveWorkbook

$wsAsset = $wb.Worksheets.Item("asset")
$wsPart  = $wb.Worksheets.Item("part")

# Rename existing headers (column E = "硬盘" -> "系统盘", column F = "RAID" -> "数据盘")
$wsPart.Range("E1").Value = "系统盘"
$wsPart.Range("F1").Value = "数据盘"

# Append the two new headers after the existing "IB卡" column (I). Inserting
# the columns (rather than just writing into J1/K1) makes the new header
# cells inherit the same bold/centered header style as the rest of row 1.
$wsPart.Columns("J:K").Insert()
$wsPart.Range("J1").Value = "模块"
$wsPart.Range("K1").Value = "配件变更"

# Make the data columns D:K a uniform wide column (matches the other wide
# columns already on the sheet, e.g. A/B/D/E/G/H which store width="20.625").
$wsPart.Range("D1:K1").EntireColumn.ColumnWidth = 19.86

# Restore the saved selections recorded in the workbook views.
$wsPart.Activate()
$wsPart.Range("A1").Select()

$wsAsset.Activate()
$wsAsset.Range("H20").Select()
